$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels were using lowercase/abbreviated forms that weren't matching
# the expected names when the data is read in and passed to the country
# lookup. Fix the header row so the column names are "Country" and
# "Country_de" (capitalized) instead of "country" / "country_de".
$ws.Range("A1").Value = "Country"
$ws.Range("B1").Value = "Country_de"
